$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1172.2222
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4836

# Row 132
$ws.Range("H132").Value = 2598699.5
$ws.Range("I132").Value = 3040716.5
$ws.Range("K132").Value = 9122149.5
$ws.Range("M132").Value = -9119619.5

# Row 135
$ws.Range("H135").Value = 4100.4326
$ws.Range("I135").Value = 694.4286
$ws.Range("J135").Value = 14696.889
$ws.Range("K135").Value = 6249.8574
$ws.Range("L135").Value = 132272.001
$ws.Range("M135").Value = -3714.8574
$ws.Range("N135").Value = -137342.001

# Row 138
$ws.Range("H138").Value = 2179.9275
$ws.Range("I138").Value = 2750
$ws.Range("J138").Value = 2007.8302
$ws.Range("K138").Value = 8250
$ws.Range("L138").Value = 6023.4906
$ws.Range("M138").Value = -3110
$ws.Range("N138").Value = -16303.4906

# Row 141
$ws.Range("H141").Value = 5991.815
$ws.Range("I141").Value = 2699.2856
$ws.Range("J141").Value = 7144.2
$ws.Range("K141").Value = 8097.8568
$ws.Range("L141").Value = 21432.6
$ws.Range("M141").Value = -2917.8568
$ws.Range("N141").Value = -31792.6

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 668.1667
$ws.Range("I74").Value = 471
$ws.Range("J74").Value = 1062.5
$ws.Range("K74").Value = 471
$ws.Range("L74").Value = 1062.5
$ws.Range("M74").Value = 403
$ws.Range("N74").Value = -2810.5

# Row 77
$ws.Range("H77").Value = 668.1667
$ws.Range("I77").Value = 471
$ws.Range("J77").Value = 1062.5
$ws.Range("K77").Value = 2355
$ws.Range("L77").Value = 5312.5
$ws.Range("M77").Value = 2013
$ws.Range("N77").Value = -14048.5

# Row 102
$ws.Range("H102").Value = 1251.8334
$ws.Range("I102").Value = 836.6667
$ws.Range("J102").Value = 1667
$ws.Range("K102").Value = 836.6667
$ws.Range("L102").Value = 1667
$ws.Range("M102").Value = 785.3333
$ws.Range("N102").Value = -4911

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2018.2354
$ws.Range("I105").Value = 1837.2727
$ws.Range("J105").Value = 2350
$ws.Range("K105").Value = 1837.2727
$ws.Range("L105").Value = 2350
$ws.Range("M105").Value = -90.27269999999999
$ws.Range("N105").Value = -5844

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3880252.8
$ws.Range("I31").Value = 3109.524
$ws.Range("J31").Value = 7581162.5
$ws.Range("K31").Value = 3109.524
$ws.Range("L31").Value = 7581162.5
$ws.Range("M31").Value = -2814.524
$ws.Range("N31").Value = -7581752.5

# Row 34
$ws.Range("H34").Value = 3880252.8
$ws.Range("I34").Value = 3109.524
$ws.Range("J34").Value = 7581162.5
$ws.Range("K34").Value = 3109.524
$ws.Range("L34").Value = 7581162.5
$ws.Range("M34").Value = -2907.524
$ws.Range("N34").Value = -7581566.5

# Row 86
$ws.Range("H86").Value = 1000000000
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -999998877
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 1000000000
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").ClearContents()

# Row 99
$ws.Range("H99").Value = 1920.2413
$ws.Range("I99").Value = 1758.4117
$ws.Range("J99").Value = 2149.5
$ws.Range("K99").Value = 1758.4117
$ws.Range("L99").Value = 2149.5
$ws.Range("M99").Value = -260.4117000000001
$ws.Range("N99").Value = -5145.5

# Row 105
$ws.Range("H105").Value = 644.4737
$ws.Range("I105").Value = 462.18182
$ws.Range("K105").Value = 462.18182
$ws.Range("M105").Value = 1284.81818

# Row 126
$ws.Range("H126").Value = 1920.2413
$ws.Range("I126").Value = 1758.4117
$ws.Range("J126").Value = 2149.5
$ws.Range("K126").Value = 5275.2351
$ws.Range("L126").Value = 6448.5
$ws.Range("M126").Value = -2805.2351
$ws.Range("N126").Value = -11388.5

# Row 132
$ws.Range("H132").Value = 3292179.2
$ws.Range("I132").Value = 2262.95
$ws.Range("J132").Value = 6947642
$ws.Range("K132").Value = 6788.849999999999
$ws.Range("L132").Value = 20842926
$ws.Range("M132").Value = -4258.849999999999
$ws.Range("N132").Value = -20847986

# Row 134
$ws.Range("H134").Value = 1137.125
$ws.Range("I134").Value = 1126.4286
$ws.Range("J134").Value = 1212
$ws.Range("K134").Value = 3379.2858
$ws.Range("L134").Value = 3636
$ws.Range("M134").Value = -844.2857999999997
$ws.Range("N134").Value = -8706

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 736.1799999999999
$ws.Range("J131").Value = 789.8652
$ws.Range("L131").Value = 2369.5956
$ws.Range("N131").Value = -12449.5956

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2954.5715
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 2697
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 2697
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -4693

# Row 83
$ws.Range("H83").Value = 2954.5715
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 2697
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 13485
$ws.Range("M83").Value = -17508
$ws.Range("N83").Value = -23469

# Row 122
$ws.Range("H122").Value = 40003096
$ws.Range("I122").Value = 100003520
$ws.Range("K122").Value = 300010560
$ws.Range("M122").Value = -300008110

# Row 126
$ws.Range("H126").Value = 2057.5881
$ws.Range("J126").Value = 1343.2222
$ws.Range("L126").Value = 4029.6666
$ws.Range("N126").Value = -8969.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2637.5
$ws.Range("I40").Value = 2216.6667
$ws.Range("J40").Value = 3900
$ws.Range("K40").Value = 2216.6667
$ws.Range("L40").Value = 3900
$ws.Range("M40").Value = -2080.6667
$ws.Range("N40").Value = -4172

# Row 68
$ws.Range("H68").Value = 20835126
$ws.Range("I68").Value = 1280
$ws.Range("J68").Value = 31252050
$ws.Range("K68").Value = 1280
$ws.Range("L68").Value = 31252050
$ws.Range("M68").Value = -531
$ws.Range("N68").Value = -31253548

# Row 71
$ws.Range("H71").Value = 20835126
$ws.Range("I71").Value = 1280
$ws.Range("J71").Value = 31252050
$ws.Range("K71").Value = 6400
$ws.Range("L71").Value = 156260250
$ws.Range("M71").Value = -2656
$ws.Range("N71").Value = -156267738

# Row 132
$ws.Range("H132").Value = 4454
$ws.Range("I132").Value = 4816.6665
$ws.Range("J132").Value = 4143.143
$ws.Range("K132").Value = 14449.9995
$ws.Range("L132").Value = 12429.429
$ws.Range("M132").Value = -11919.9995
$ws.Range("N132").Value = -17489.429

# Row 136
$ws.Range("H136").Value = 10435.462
$ws.Range("I136").Value = 11605.546
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 34816.638
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -32266.638
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1540.2
$ws.Range("I136").Value = 1546.7693
$ws.Range("J136").Value = 1497.5
$ws.Range("K136").Value = 4640.3079
$ws.Range("L136").Value = 4492.5
$ws.Range("M136").Value = -2090.3079
$ws.Range("N136").Value = -9592.5

# Row 138
$ws.Range("H138").Value = 42328.6
$ws.Range("J138").Value = 42328.6
$ws.Range("L138").Value = 42328.6
$ws.Range("N138").Value = -52608.6

Write-Output "Applied Garuda_Profits updates"
